$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8164013333333333
$ws.Range("H2").Value = 2.449204
$ws.Range("I2").Value = 0.05618115571687973
$ws.Range("J2").Value = 0.05618115571687973
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 9.100299936861335
$ws.Range("R2").Value = 81.90269943175201
$ws.Range("S2").Value = 0.01457792047542534
$ws.Range("T2").Value = 0.01457792047542534
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8164013333333333
$ws.Range("H3").Value = 2.449204
$ws.Range("I3").Value = 0.05618115571687973
$ws.Range("J3").Value = 0.05618115571687973
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 22.602476664796
$ws.Range("R3").Value = 203.422289983164
$ws.Range("S3").Value = 0.03620727994166482
$ws.Range("T3").Value = 0.03620727994166482
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8164013333333333
$ws.Range("H4").Value = 2.449204
$ws.Range("I4").Value = 0.05618115571687973
$ws.Range("J4").Value = 0.05618115571687973
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 3.368437340343556
$ws.Range("R4").Value = 30.315936063092
$ws.Range("S4").Value = 0.005395955299789569
$ws.Range("T4").Value = 0.005395955299789569
$ws.Range("I5").Value = 0.8862323361798529
$ws.Range("J5").Value = 0.8862323361798529
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 143.553117946252
$ws.Range("R5").Value = 1291.978061516268
$ws.Range("S5").Value = 0.2299601059238916
$ws.Range("T5").Value = 0.2299601059238916
$ws.Range("I6").Value = 0.8862323361798529
$ws.Range("J6").Value = 0.8862323361798529
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.5711534745053068
$ws.Range("T6").Value = 0.5711534745053068
$ws.Range("I7").Value = 0.8862323361798529
$ws.Range("J7").Value = 0.8862323361798529
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.08511875575065442
$ws.Range("T7").Value = 0.08511875575065442
$ws.Range("G8").Value = 0.8368233333333334
$ws.Range("I8").Value = 0.05758650810326746
$ws.Range("J8").Value = 0.05758650810326746
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 9.327940825873336
$ws.Range("R8").Value = 83.95146743286001
$ws.Range("S8").Value = 0.01494258216789661
$ws.Range("T8").Value = 0.01494258216789661
$ws.Range("G9").Value = 0.8368233333333334
$ws.Range("I9").Value = 0.05758650810326746
$ws.Range("J9").Value = 0.05758650810326746
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.03711299266012603
$ws.Range("T9").Value = 0.03711299266012603
$ws.Range("G10").Value = 0.8368233333333334
$ws.Range("I10").Value = 0.05758650810326746
$ws.Range("J10").Value = 0.05758650810326746
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("Q10").Value = 3.452697647812223
$ws.Range("S10").Value = 0.005530933275244823
$ws.Range("T10").Value = 0.005530933275244823
